$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three changed labels in column B
$ws.Range("B42").Value = "payment_det"
$ws.Range("B14").Value = "f_age"
$ws.Range("B43").Value = "payment"

# Update the view: top-left cell and selection
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A26").Select()
